# Add a new "QUESTION & ANSWER" section (row 20) to the API reference sheet,
# with its action, endpoint, sample JSON response and extra note - mirroring
# the layout used by the other section headers (SECTIONS, CHAPTERS, ACTIVITY).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionJson = @'
{
        "number": 1,
        "id": 1,
        "q_type": "IDENT",
        "question_name": "This is a question #1",
        "choices": []
    },
    {
        "number": 2,
        "id": 2,
        "q_type": "MULT",
        "question_name": "This is a question multiple #2",
        "choices": [
            {
                "description": "A. sample"
            },
            {
                "description": "B. sample"
            }
        ]
    }
'@

# --- New row 20: QUESTION & ANSWER section header row ---
$ws.Range("A20").Value = "QUESTION & ANSWER"
$ws.Range("B20").Value = "list view"
$ws.Range("C20").Value = "/api/question/list/<activity>/"
$ws.Range("D20").Value = $questionJson
$ws.Range("E20").Value = "get only method"

# Formatting to match the other section-header rows (e.g. row 11 "ACTIVITY"):
# centered, non-wrapping text for the label/action/url/extra columns ...
$headerCells = @("A20", "B20", "C20", "E20")
foreach ($addr in $headerCells) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4108     # xlCenter
    $c.WrapText = $false
}

# ... and a wrapping, general/bottom aligned JSON sample cell, as used for
# the other sample-response cells (D3, D6, D10, D13, D14, D17, D18, ...).
$dCell = $ws.Range("D20")
$dCell.HorizontalAlignment = 1     # xlGeneral
$dCell.VerticalAlignment = -4107  # xlBottom
$dCell.WrapText = $true

# The row holds a long multi-line JSON sample, so it needs extra height,
# same as it is given in the workbook.
$ws.Rows.Item(20).RowHeight = 236.55

# Scroll the view down and select the new cell, like the author did when
# finishing the edit.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D25").Select()
